$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they pick up the same bold/border/center style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-28
$iValues = @(6,6,6,2,7,1,7,7,5,1,1,1,1,1,1,1,1,1,1,7,1,11,1,1,1,5,3)
$jValues = @(8,6,6,6,7,5,8,8,7,5,4,5,6,5,4,6,4,4,5,9,5,11,5,5,3,7,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

$ws.Range("A1").Select()
